# Daily attendance processing - 2026-01-26 12:59:24
# In the "Recorded By" column, swap the order of the two recorder names so
# that "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1
$firstCol = $usedRange.Column
$lastCol = $usedRange.Columns.Count + $firstCol - 1

# Locate the "Recorded By" column dynamically (falls back to column 7 / G).
$recordedByCol = 7
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
